# Applies updated cryptocurrency price/volume figures scraped on 2023-01-06.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new literal text value (Price / Volume(1h) columns).
$updates = [ordered]@{
    "D2" = "258.87"
    "E2" = "0.43%"
    "D3" = "26.93"
    "E3" = "-0.88%"
    "D4" = "4.685"
    "E4" = "0.42%"
    "D5" = "0.06036"
    "E5" = "2.51%"
    "D6" = "6.680"
    "E6" = "0.55%"
    "D7" = "0.8587"
    "E7" = "0.10%"
    "D8" = "0.9278"
    "E8" = "-1.69%"
    "D9" = "0.1395"
    "E9" = "-0.91%"
    "D10" = "0.04828"
    "E10" = "25.55%"
    "D11" = "0.07064"
    "E11" = "-0.24%"
    "D12" = "0.03128"
    "E12" = "-1.58%"
    "D13" = "0.09130"
    "E13" = "-0.55%"
    "D14" = "0.001533"
    "E14" = "-0.53%"
    "D15" = "0.0006067"
    "E15" = "0.66%"
    "D16" = "0.006013"
    "E16" = "-2.87%"
    "D18" = "3.166"
    "E18" = "-1.26%"
    "E19" = "-1.75%"
    "E20" = "0.41%"
    "D21" = "0.1297"
    "E21" = "0.36%"
    "D22" = "4.122"
    "E22" = "6.23%"
    "D23" = "0.04241"
    "E23" = "0.79%"
    "E24" = "-0.33%"
    "D25" = "0.004037"
    "E25" = "-6.08%"
    "D26" = "0.0001200"
    "E26" = "0.02%"
    "D27" = "0.0001523"
    "E27" = "-21.33%"
    "D40" = "0.03841"
    "E40" = "0.16%"
    "D41" = "0.1114"
    "E41" = "1.01%"
    "E42" = "-37.24%"
    "D43" = "0.01527"
    "E43" = "33.27%"
    "D44" = "0.002200"
    "E44" = "0.03%"
    "D45" = "0.00005131"
    "E45" = "-5.96%"
    "E46" = "0.03%"
    "D48" = "0.1353"
    "E48" = "3.87%"
    "D49" = "0.00002100"
    "E49" = "0.03%"
    "D50" = "0.0002000"
    "E50" = "0.03%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage so values like "258.87" / "0.43%" are not
    # auto-converted to numbers/percentages by Excel.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    # Reset the style back to Normal so no stray number-format styling
    # is left behind on the cell (matches original "General" formatting).
    $cell.Style = "Normal"
}
